$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking values
# (e.g. "0.9992", "318.90") are stored as text, matching the source data,
# then restore the default "Normal" style so no stray style index remains.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.399.25"
$ws.Range("D3").Value = "1.829.99"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").Value = "318.90"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "0.5330"
$ws.Range("E7").Value = "  -1.24%  "
$ws.Range("D8").Value = "0.3990"
$ws.Range("E8").Value = "  +5.54%  "
$ws.Range("D9").Value = "0.07575"
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("D10").Value = "41.74"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "6.320"
$ws.Range("E12").Value = "  +2.82%  "
$ws.Range("D13").Value = "7.621"
$ws.Range("E13").Value = "  +4.26%  "
$ws.Range("D14").Value = "0.9988"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "20.74"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").Value = "1.821.52"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("D17").Value = "89.90"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "0.00001073"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D19").Value = "0.06591"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").Value = "17.57"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").Value = "0.9983"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").Value = "6.064"
$ws.Range("E22").Value = "  +2.15%  "
$ws.Range("D23").Value = "28.404.75"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("D25").Value = "2.105"
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("D26").Value = "156.69"
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").Value = "2.419"
$ws.Range("E28").Value = "  +4.07%  "
$ws.Range("D29").Value = "2.032.30"
$ws.Range("E29").Value = "  +1.64%  "
$ws.Range("D30").Value = "123.71"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").Value = "1.113"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").Value = "0.1102"
$ws.Range("E32").Value = "  +4.56%  "
$ws.Range("D33").Value = "3.689"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("D34").Value = "5.610"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").Value = "0.07369"
$ws.Range("E35").Value = "  +12.92%  "
$ws.Range("D36").Value = "0.2252"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").Value = "5.253"
$ws.Range("E37").Value = "  +4.42%  "
$ws.Range("D38").Value = "0.02332"
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("D39").Value = "8.798"
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("D40").Value = "11.35"
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("D41").Value = "0.6265"
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("D42").Value = "1.196"
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("D43").Value = "1.416"
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("D44").Value = "13.44"
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("D45").Value = "3.710"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("D46").Value = "0.5825"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "125.73"
$ws.Range("D48").Value = "1.971"
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("D49").Value = "1.193"
$ws.Range("E49").Value = "  -1.65%  "
$ws.Range("D50").Value = "0.06885"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "71.55"
$ws.Range("E51").Value = "  -1.05%  "

$priceRange.Style = "Normal"
